# Add the two new "RAM" result sheets (taskset creation / optimization-algo
# benchmark data) after the existing ro_FLASH-code_CCM sheet, mirroring the
# layout/style of the pre-existing ro_* sheets.

$wb = $excel.ActiveWorkbook

# Sheet used as the formatting template for the new sheets (last existing
# sheet, already carries the bold/bordered/centered header style).
$template = $wb.Worksheets.Item("ro_FLASH-code_CCM")

# ---------------------------------------------------------------------
# Sheet: ro_RAM-code_FLASH
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$s5 = $wb.Worksheets.Add($null, $lastSheet)
$s5.Name = "ro_RAM-code_FLASH"

$s5.Range("B1").Value = "'24"
$s5.Range("C1").Value = "'48"
$s5.Range("D1").Value = "'72"

$s5.Range("A2").Value = "intensity"
$s5.Range("B2").Value = 14992
$s5.Range("C2").Value = 26643
$s5.Range("D2").Value = 33958

$s5.Range("A3").Value = "runtime"
$s5.Range("B3").Value = 0.33726
$s5.Range("C3").Value = 0.18166
$s5.Range("D3").Value = 0.15584

$s5.Range("A4").Value = "timestamp"
$s5.Range("B4").Value = "(1581.52, 1918.78)"
$s5.Range("C4").Value = "(3427.48, 3609.14)"
$s5.Range("D4").Value = "(4949.36, 5105.2)"

$s5.Range("A5").Value = "energy"
$s5.Range("B5").Value = 16.685
$s5.Range("C5").Value = 15.972
$s5.Range("D5").Value = 17.464

# Copy the header/label formatting (bold, thin border, centered) from the
# template sheet without disturbing the text we just wrote.
$template.Range("B1:D5").Copy()
$s5.Range("B1:D5").PasteSpecial(-4122)  # xlPasteFormats
$template.Range("A2:A5").Copy()
$s5.Range("A2:A5").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------
# Sheet: ro_RAM-code_CCM
# ---------------------------------------------------------------------
$s6 = $wb.Worksheets.Add($null, $s5)
$s6.Name = "ro_RAM-code_CCM"

$s6.Range("B1").Value = "'24"
$s6.Range("C1").Value = "'48"
$s6.Range("D1").Value = "'72"

$s6.Range("A2").Value = "intensity"
$s6.Range("B2").Value = 12297
$s6.Range("C2").Value = 24223
$s6.Range("D2").Value = 36222

$s6.Range("A3").Value = "runtime"
$s6.Range("B3").Value = 0.33694
$s6.Range("C3").Value = 0.16872
$s6.Range("D3").Value = 0.1128

$s6.Range("A4").Value = "timestamp"
$s6.Range("B4").Value = "(2114.38, 2451.32)"
$s6.Range("C4").Value = "(3804.48, 3973.2)"
$s6.Range("D4").Value = "(5300.54, 5413.34)"

$s6.Range("A5").Value = "energy"
$s6.Range("B5").Value = 13.673
$s6.Range("C5").Value = 13.487
$s6.Range("D5").Value = 13.483

$template.Range("B1:D5").Copy()
$s6.Range("B1:D5").PasteSpecial(-4122)  # xlPasteFormats
$template.Range("A2:A5").Copy()
$s6.Range("A2:A5").PasteSpecial(-4122)  # xlPasteFormats

Write-Host "Added sheets ro_RAM-code_FLASH and ro_RAM-code_CCM"
